$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old rows 6 and 7) that no longer exist in the new data
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# Row 2
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Il2"
$ws.Range("C2").Value = "Il2rb"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07963266666666667
$ws.Range("H2").Value = 0.238898
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.805089
$ws.Range("N2").Value = 2.415267
$ws.Range("O2").Value = 0.4118548214512568
$ws.Range("P2").Value = 0.4156154679278413
$ws.Range("Q2").Value = 0.06411138397400001
$ws.Range("R2").Value = 0.577002455766
$ws.Range("S2").Value = 0.4118548214512568
$ws.Range("T2").Value = 0.4156154679278413

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Il2"
$ws.Range("C3").Value = "Il2rb"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07963266666666667
$ws.Range("H3").Value = 0.238898
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9511543333333333
$ws.Range("N3").Value = 2.853463
$ws.Range("O3").Value = 0.4865766370271973
$ws.Range("P3").Value = 0.4910195684202955
$ws.Range("Q3").Value = 0.07574295597488889
$ws.Range("R3").Value = 0.681686603774
$ws.Range("S3").Value = 0.4865766370271973
$ws.Range("T3").Value = 0.4910195684202955

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Il2"
$ws.Range("C4").Value = "Il2rb"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.07963266666666667
$ws.Range("H4").Value = 0.238898
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.053063
$ws.Range("N4").Value = 0.106126
$ws.Range("O4").Value = 0.0271451384762033
$ws.Range("P4").Value = 0.0182620004948977
$ws.Range("Q4").Value = 0.004225548191333334
$ws.Range("R4").Value = 0.025353289148
$ws.Range("S4").Value = 0.0271451384762033
$ws.Range("T4").Value = 0.0182620004948977

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Il2"
$ws.Range("C5").Value = "Il2rb"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.07963266666666667
$ws.Range("H5").Value = 0.238898
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.145482
$ws.Range("N5").Value = 0.436446
$ws.Range("O5").Value = 0.07442340304534249
$ws.Range("P5").Value = 0.07510296315696552
$ws.Range("Q5").Value = 0.011585119612
$ws.Range("R5").Value = 0.104266076508
$ws.Range("S5").Value = 0.07442340304534249
$ws.Range("T5").Value = 0.07510296315696552
